# Re-sort the historical-distance rows (time bucket analysis).
# The underlying data set (title / timestamp / historical distance /
# time bucket / uri) is unchanged per article, only the row order in
# which the articles are listed changes (as if the generating script
# had been re-run after adding a new JSON source to sort by).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 15

# Snapshot the current rows (A..E) before touching anything.
$data = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @(
        $ws.Cells.Item($r, 1).Value(),
        $ws.Cells.Item($r, 2).Value(),
        $ws.Cells.Item($r, 3).Value(),
        $ws.Cells.Item($r, 4).Value(),
        $ws.Cells.Item($r, 5).Value()
    )
    $data += ,$row
}

# New row order, expressed as 0-based indices into $data
# (index 0 == old row 2, index 13 == old row 15).
$order = @(1, 4, 0, 3, 2, 7, 5, 6, 12, 11, 13, 8, 9, 10)

# Drop the existing hyperlinks; they will be re-created against the
# (possibly different) uri that now lands on each row.
$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $order.Count; $i++) {
    $src = $data[$order[$i]]
    $r = $firstRow + $i

    $ws.Cells.Item($r, 1).Value = $src[0]
    $ws.Cells.Item($r, 2).Value = $src[1]
    $ws.Cells.Item($r, 3).Value = $src[2]
    $ws.Cells.Item($r, 4).Value = $src[3]
    $ws.Cells.Item($r, 5).Value = $src[4]

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 5), $src[4])
}
